$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New unified StatQuery text (replaces the old per-tab StatQuery formula in column C)
$newStatQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['Cavalier King Charles Spaniel']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Column width adjustments (column B keeps its existing width; C and D get new, narrower
# best-fit widths now that column C holds the unified StatQuery text and D holds the
# (shorter) dbExcel filename values)
$ws.Columns.Item(3).ColumnWidth = 60.5
$ws.Columns.Item(4).ColumnWidth = 49.25

# View changes: zoom reset to 100% (Normal view) and selection moved to B4
$excel.ActiveWindow.Zoom = 100
$ws.Range("B4").Select()
